$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.831.04"
$ws.Range("E2").Value = "  +4.15%  "
$ws.Range("D3").Value = "2.772.22"
$ws.Range("E3").Value = "  +4.53%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.65"
$ws.Range("E5").Value = "  +4.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "115.22"
$ws.Range("E6").Value = "  +1.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.547"
$ws.Range("E7").Value = "  +4.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.574"
$ws.Range("E9").Value = "  +4.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.02"
$ws.Range("E10").Value = "  +5.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0854"
$ws.Range("E11").Value = "  +4.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.92"
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.61"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").Value = "3.216.37"
$ws.Range("E15").Value = "  +4.89%  "
$ws.Range("D16").Value = "2.759.63"
$ws.Range("E16").Value = "  +4.07%  "
$ws.Range("D17").Value = "51.812.15"
$ws.Range("E17").Value = "  +4.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.877"
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.20"
$ws.Range("E19").Value = "  +8.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.02"
$ws.Range("E20").Value = "  +4.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.17"
$ws.Range("E21").Value = "  -1.63%  "
$ws.Range("D22").Value = "0.0₃0977"
$ws.Range("E22").Value = "  +2.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "275.55"
$ws.Range("E23").Value = "  +2.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.98"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.75"
$ws.Range("E25").Value = "  +6.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.50"
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.16"
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.69"
$ws.Range("E31").Value = "  -0.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.01"
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.68"
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0816"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.09"
$ws.Range("E36").Value = "  +2.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.90"
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.21"
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0382"
$ws.Range("E40").Value = "  +10.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.65"
$ws.Range("E41").Value = "  +23.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.34"
$ws.Range("E42").Value = "  +2.56%  "
$ws.Range("E43").Value = "  +3.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.84"
$ws.Range("E44").Value = "  -1.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.25"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("D46").Value = "2.065.62"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.30"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("E48").Value = "  +1.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.54"
$ws.Range("E49").Value = "  +4.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.890"
$ws.Range("E50").Value = "  +13.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.86"
$ws.Range("E51").Value = "  -0.99%  "
